# Swap the data contents of rows 2 and 3 on the "Artfynd" sheet for the
# columns that actually differ between the two rows (A, B, E, F, G, H, Q, R).
# All other columns are identical between row 2 and row 3, so only these
# need to be touched to realize the row swap described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $columns) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $val3
    $ws.Range($addr3).Value = $val2
}
